$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New purchase entries (rows 26-31), continuing the existing cost-tracking table.
# Columns: A=Provider, B=Description, C=unit cost ($), D=Person responsible
$ws.Range("A26").Value = "eBay"
$ws.Range("B26").Value = "Radio Chino"
$ws.Range("C26").Value = 21.8
$ws.Range("D26").Value = "David"

$ws.Range("A27").Value = "Amazon"
$ws.Range("B27").Value = "Cama s100 Canon"
$ws.Range("C27").Value = 156.6
$ws.Range("D27").Value = "David"

$ws.Range("A28").Value = "3DR"
$ws.Range("B28").Value = "Ardupilot"
$ws.Range("C28").Value = 243
$ws.Range("D28").Value = "Julio"

$ws.Range("A29").Value = "Canon"
$ws.Range("B29").Value = "Camara ELPH s130"
$ws.Range("C29").Value = 100
$ws.Range("D29").Value = "Julio"

$ws.Range("A30").Value = "Hobby King"
$ws.Range("B30").Value = "Radio Spektrum"
$ws.Range("C30").Value = 166
$ws.Range("D30").Value = "Julio"

$ws.Range("A31").Value = "Hobby Town"
$ws.Range("B31").Value = "Speed Controller"
$ws.Range("C31").Value = 55
$ws.Range("D31").Value = "Julio"

# Restore the cursor/selection to match where the author left off editing.
$ws.Range("G27").Select()
